$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$clothingValues = @{
    2 = 'Jumpsuit,Dress'
    3 = 'Jumpsuit,Tee'
    4 = 'Halter,Parka'
    5 = 'Halter,Jumpsuit'
    6 = 'Jumpsuit,Blouse'
    7 = 'Jumpsuit,Parka'
    8 = 'Jumpsuit,Dress'
    9 = 'Parka,Dress'
    10 = 'Jumpsuit,Cutoffs'
    11 = 'Jumpsuit,Dress'
    12 = 'Jumpsuit,Halter'
    13 = 'Jumpsuit,Blouse'
    14 = 'Tee,Kaftan'
    15 = 'Kaftan,Tee'
    16 = 'Jumpsuit,Tee'
    17 = 'Dress,Trunks'
    18 = 'Kaftan,Trunks'
    19 = 'Jumpsuit,Dress'
    20 = 'Parka,Gauchos'
    21 = 'Parka,Halter'
    22 = 'Jumpsuit,Caftan'
    23 = 'Tee,Jumpsuit'
    24 = 'Jumpsuit,Blouse'
    25 = 'Jumpsuit,Dress'
    26 = 'Jumpsuit,Blazer'
    27 = 'Jumpsuit,Kaftan'
    28 = 'Jumpsuit,Blouse'
    29 = 'Blazer,Jumpsuit'
    30 = 'Jumpsuit,Dress'
    31 = 'Jumpsuit,Kaftan'
    32 = 'Jumpsuit,Trunks'
    33 = 'Jumpsuit,Dress'
    34 = 'Parka,Kaftan'
    35 = 'Parka,Jumpsuit'
    36 = 'Parka,Tee'
    37 = 'Jumpsuit,Dress'
    38 = 'Tee,Jumpsuit'
    39 = 'Sweatpants,Jumpsuit'
    40 = 'Jumpsuit,Dress'
    41 = 'Blouse,Jumpsuit'
    42 = 'Halter,Blazer'
    43 = 'Blouse,Jumpsuit'
    44 = 'Jumpsuit,Halter'
    45 = 'Jumpsuit,Kaftan'
    46 = 'Jumpsuit,Kaftan'
    47 = 'Tee,Kaftan'
    48 = 'Jumpsuit,Parka'
    49 = 'Parka,Halter'
    50 = 'Kaftan,Parka'
    51 = 'Kaftan,Parka'
    52 = 'Jumpsuit,Kaftan'
    53 = 'Jumpsuit,Halter'
    54 = 'Jumpsuit,Halter'
    55 = 'Halter,Trunks'
    56 = 'Jumpsuit,Halter'
    57 = 'Halter,Tee'
    58 = 'Halter,Jumpsuit'
    59 = 'Jumpsuit,Halter'
    60 = 'Jumpsuit,Blouse'
    61 = 'Jumpsuit,Parka'
    62 = 'Jumpsuit,Halter'
    63 = 'Kaftan,Tee'
    64 = 'Halter,Jumpsuit'
    65 = 'Jumpsuit,Blouse'
    66 = 'Jumpsuit,Halter'
    67 = 'Halter,Blazer'
    68 = 'Jumpsuit,Blouse'
    69 = 'Blazer,Halter'
    70 = 'Jumpsuit,Blouse'
    71 = 'Kaftan,Tee'
    72 = 'Jumpsuit,Blouse'
    73 = 'Kaftan,Tee'
    74 = 'Jumpsuit,Kaftan'
    75 = 'Parka,Halter'
    76 = 'Blouse,Jumpsuit'
    77 = 'Jumpsuit,Halter'
    78 = 'Blazer,Jumpsuit'
    79 = 'Kaftan,Jumpsuit'
    80 = 'Blouse,Halter'
    81 = 'Jumpsuit,Dress'
    82 = 'Kaftan,Tee'
    83 = 'Halter,Blazer'
    84 = 'Jumpsuit,Halter'
    85 = 'Parka,Caftan'
    86 = 'Blazer,Parka'
    87 = 'Top,Jumpsuit'
    88 = 'Parka,Blouse'
}

foreach ($row in $clothingValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $clothingValues[$row]
}

$wb.Save()
